{"js": "// The edit:\n//  - The paragraph \"I was previously a [bookmark]student from Guangdong\n//    Technion. While I terminated my study due to some personal reasons,\n//    currently, I have several papers under review with PNAS, and I am in\n//    the process of preparing other manuscripts for submission to PNAS. \"\n//    becomes a single run reading \"I was previously an undergraduate\n//    student from Guangdong Technion. While my study was terminated due to\n//    some personal reasons, currently, I have several papers under review\n//    with PNAS, and I am working with professors from Berkeley in the\n//    process of preparing other manuscripts for submission to PNAS. \"\n//  - The \"_GoBack\" bookmark that used to sit mid-sentence moves down into\n//    the (already existing) blank paragraph right after this one.\n\nconst body = context.document.body;\n\n// Locate the target paragraph robustly by searching for a distinctive,\n// stable substring that survives the edit (rather than hard-coding an\n// index), then resolve it to its containing Paragraph object.\nconst results = body.search(\"Guangdong Technion\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the target paragraph ('Guangdong Technion').\");\n}\n\nconst hitParagraphs = results.items[0].paragraphs;\nhitParagraphs.load(\"items\");\nawait context.sync();\n\nconst targetParagraph = hitParagraphs.items[0];\n\n// Replace the whole paragraph's text (which also removes the old inline\n// bookmark, since it gets overwritten along with the rest of the content)\n// with the revised wording, as a single run.\nconst newText =\n  \"I was previously an undergraduate student from Guangdong Technion. \" +\n  \"While my study was terminated due to some personal reasons, currently, \" +\n  \"I have several papers under review with PNAS, and I am working with \" +\n  \"professors from Berkeley in the process of preparing other manuscripts \" +\n  \"for submission to PNAS. \";\n\ntargetParagraph.getRange(\"Whole\").insertText(newText, \"Replace\");\nawait context.sync();\n\n// Re-insert the \"_GoBack\" bookmark into the blank paragraph that follows\n// (the document already had an empty paragraph there).\nconst nextParagraph = targetParagraph.getNext();\nnextParagraph.getRange(\"Whole\").insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The edit:\n#  - The paragraph \"I was previously a [bookmark]student from Guangdong\n#    Technion. While I terminated my study due to some personal reasons,\n#    currently, I have several papers under review with PNAS, and I am in\n#    the process of preparing other manuscripts for submission to PNAS. \"\n#    becomes a single run reading \"I was previously an undergraduate\n#    student from Guangdong Technion. While my study was terminated due to\n#    some personal reasons, currently, I have several papers under review\n#    with PNAS, and I am working with professors from Berkeley in the\n#    process of preparing other manuscripts for submission to PNAS. \"\n#  - The \"_GoBack\" bookmark that used to sit mid-sentence moves down into\n#    the (already existing) blank paragraph right after this one.\n\n$d = $word.ActiveDocument\n\n# Locate the target paragraph robustly via Find (rather than a hard-coded\n# paragraph index) using a distinctive, stable substring.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$found = $find.Execute(\"Guangdong Technion\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\nif (-not $found) {\n    throw \"Could not find the target paragraph ('Guangdong Technion').\"\n}\n\n# Resolve the full containing paragraph from the (now collapsed) find hit.\n$targetPara = $rng.Paragraphs(1)\n\n# Replace the whole paragraph's text (excluding the trailing paragraph\n# mark) with the revised wording. Assigning .Text on the whole-paragraph\n# range collapses all runs into one, which also removes the old inline\n# \"_GoBack\" bookmark along with the rest of the overwritten content.\n$paraRng = $targetPara.Range\n$paraRng.MoveEnd(1, -1) | Out-Null\n$paraRng.Text = \"I was previously an undergraduate student from Guangdong Technion. While my study was terminated due to some personal reasons, currently, I have several papers under review with PNAS, and I am working with professors from Berkeley in the process of preparing other manuscripts for submission to PNAS. \"\n\n# Re-insert the \"_GoBack\" bookmark into the blank paragraph that follows\n# (the document already had an empty paragraph there).\n$nextPara = $targetPara.Next()\n$d.Bookmarks.Add(\"_GoBack\", $nextPara.Range)\n"}
